# Gear change implementation in getWheelTorques3 function and results from portimao
#
# - Engine sheet: final drive / gear ratio bumped from 3.1 to 3.5, and a new
#   "gear_ratio" lookup column (pointing at GearRatio.xlsx) is added next to
#   the existing "map" column.
# - Susp sheet: updated result value from the portimao run (0.51 -> 0.56).
# - The Brake sheet becomes the active/selected sheet (was Engine before).
# - Misc selection/active-cell bookkeeping left behind by the interactive
#   session is replayed on each sheet.

$wb = $excel.ActiveWorkbook

# --- Engine sheet: gear ratio change + new gear_ratio/GearRatio.xlsx column ---
$wsEngine = $wb.Worksheets.Item("Engine")
$wsEngine.Range("A2").Value = 3.5
$wsEngine.Range("E1").Value = "gear_ratio"
$wsEngine.Range("E2").Value = "GearRatio.xlsx"

# --- Susp sheet: updated value from the portimao results ---
$wsSusp = $wb.Worksheets.Item("Susp")
$wsSusp.Range("A2").Value = 0.56000000000000005

# --- Replay selections / active cells left over from the editing session ---
$wsMassInertia = $wb.Worksheets.Item("MassInertia")
$wsMassInertia.Range("G21").Select() | Out-Null

$wsAero = $wb.Worksheets.Item("Aero")
$wsAero.Range("B3").Select() | Out-Null

$wsEngine.Range("E3").Select() | Out-Null

# Brake becomes the active sheet (tabSelected moves off Engine and onto Brake)
$wsBrake = $wb.Worksheets.Item("Brake")
$wsBrake.Activate() | Out-Null
$wsBrake.Range("J23").Select() | Out-Null
